# TENNESSEE_2017.xlsx data-cleaning edit
#
# 1. Rename the header row to the short machine-friendly column names.
# 2. Title-case the Spanish connector words (de/del/la/las/los/el/y) that
#    appear lower-cased inside place names in columns A and B.
# 3. Correct one floating point rounding artifact in D574.
# 4. Drop the trailing metadata/footnote rows (1349-1354) and shrink the
#    sheet's used range back down to A1:D1348.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row -> short column codes -------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case the Spanish prepositions/articles in place names ----
# (columns A = state, B = municipality; row 1 is the header we just set,
# and the data we care about runs through row 1348)
$lastRow = 1348
for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("A", "B")) {
        $cell = $ws.Range("$col$r")
        $v = $cell.Value()
        if ($v -ne $null -and $v -is [string]) {
            $nv = $v -replace '\bde\b', 'De' `
                      -replace '\bdel\b', 'Del' `
                      -replace '\blas\b', 'Las' `
                      -replace '\blos\b', 'Los' `
                      -replace '\bla\b', 'La' `
                      -replace '\bel\b', 'El' `
                      -replace '\by\b', 'Y'
            if (-not $nv.Equals($v)) {
                $cell.Value = $nv
            }
        }
    }
}

# --- 3. Fix the rounding artifact on D574 -------------------------------
$ws.Range("D574").Value = 0.00994229915667998

# --- 4. Remove the trailing metadata rows and shrink the used range ----
$ws.Range("A1349:A1354").EntireRow.Delete()
